$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 9; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 10; I = "ba"; J = "Appreciation" },
    @{ Row = 13; I = "b"; J = "Acknowledge (Backchannel)" },
    @{ Row = 18; I = "ba"; J = "Appreciation" },
    @{ Row = 19; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 34; I = "aa"; J = "Agree/Accept" },
    @{ Row = 35; I = "ba"; J = "Appreciation" },
    @{ Row = 37; I = "ba"; J = "Appreciation" },
    @{ Row = 46; I = "ba"; J = "Appreciation" },
    @{ Row = 53; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 54; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 56; I = "aa"; J = "Agree/Accept" },
    @{ Row = 58; I = "ba"; J = "Appreciation" },
    @{ Row = 60; I = "ba"; J = "Appreciation" },
    @{ Row = 77; I = "ba"; J = "Appreciation" },
    @{ Row = 79; I = "ba"; J = "Appreciation" },
    @{ Row = 84; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 85; I = "aa"; J = "Agree/Accept" },
    @{ Row = 90; I = "sv"; J = "Statement-opinion" },
    @{ Row = 95; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 96; I = "%"; J = "Uninterpretable" },
    @{ Row = 98; I = "b"; J = "Acknowledge (Backchannel)" },
    @{ Row = 101; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 102; I = "sv"; J = "Statement-opinion" },
    @{ Row = 103; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 124; I = "ba"; J = "Appreciation" },
    @{ Row = 127; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 137; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 147; I = "ba"; J = "Appreciation" },
    @{ Row = 149; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 153; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 163; I = "ba"; J = "Appreciation" },
    @{ Row = 164; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 169; I = "aa"; J = "Agree/Accept" },
    @{ Row = 179; I = "ba"; J = "Appreciation" },
    @{ Row = 191; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 192; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 199; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 201; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 206; I = "ba"; J = "Appreciation" },
    @{ Row = 212; I = "aa"; J = "Agree/Accept" },
    @{ Row = 213; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 224; I = "ba"; J = "Appreciation" },
    @{ Row = 230; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 234; I = "sv"; J = "Statement-opinion" },
    @{ Row = 241; I = "%"; J = "Uninterpretable" },
    @{ Row = 244; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 253; I = "ba"; J = "Appreciation" },
    @{ Row = 259; I = "ba"; J = "Appreciation" },
    @{ Row = 268; I = "aa"; J = "Agree/Accept" },
    @{ Row = 269; I = "aa"; J = "Agree/Accept" },
    @{ Row = 270; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 277; I = "aa"; J = "Agree/Accept" },
    @{ Row = 288; I = "sv"; J = "Statement-opinion" },
    @{ Row = 290; I = "aa"; J = "Agree/Accept" },
    @{ Row = 291; I = "%"; J = "Uninterpretable" },
    @{ Row = 299; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 300; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 301; I = "sv"; J = "Statement-opinion" },
    @{ Row = 303; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 304; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 305; I = "aa"; J = "Agree/Accept" },
    @{ Row = 308; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 309; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 310; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 312; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 316; I = "ba"; J = "Appreciation" },
    @{ Row = 324; I = "ba"; J = "Appreciation" },
    @{ Row = 340; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 342; I = "ba"; J = "Appreciation" },
    @{ Row = 351; I = "ba"; J = "Appreciation" },
    @{ Row = 354; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 359; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 367; I = "ba"; J = "Appreciation" },
    @{ Row = 370; I = "sv"; J = "Statement-opinion" },
    @{ Row = 371; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 388; I = "aa"; J = "Agree/Accept" },
    @{ Row = 389; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 395; I = "sv"; J = "Statement-opinion" },
    @{ Row = 404; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 407; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 414; I = "ba"; J = "Appreciation" },
    @{ Row = 417; I = "%"; J = "Uninterpretable" },
    @{ Row = 421; I = "sv"; J = "Statement-opinion" },
    @{ Row = 422; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 428; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 443; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 450; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 454; I = "aa"; J = "Agree/Accept" },
    @{ Row = 472; I = "ba"; J = "Appreciation" },
    @{ Row = 473; I = "sv"; J = "Statement-opinion" },
    @{ Row = 476; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

